$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.25682061283031
$ws.Range("C2").Value = 9.707821550004406
$ws.Range("E2").Value = 16.6088673869258
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.584385311194781
$ws.Range("I2").Value = 15.87132767322489
$ws.Range("N2").Value = 15.49445252527274
$ws.Range("O2").Value = 16.67332912509209
$ws.Range("B3").Value = 11.59573960824411
$ws.Range("C3").Value = 9.197241677216487
$ws.Range("E3").Value = 15.65914554091905
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.586552551079531
$ws.Range("I3").Value = 15.97839623694138
$ws.Range("N3").Value = 15.52777245959841
$ws.Range("O3").Value = 16.69704255660327
$ws.Range("B4").Value = 11.17056138294702
$ws.Range("C4").Value = 8.8671046754313
$ws.Range("E4").Value = 15.05055784003633
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.587951646056756
$ws.Range("I4").Value = 16.04889509649628
$ws.Range("N4").Value = 15.55004499131614
$ws.Range("O4").Value = 16.71781474913521
$ws.Range("B5").Value = 10.99261658633563
$ws.Range("C5").Value = 8.728460268027954
$ws.Range("E5").Value = 14.7964190552931
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.588539048831045
$ws.Range("I5").Value = 16.07881571445554
$ws.Range("N5").Value = 15.55957806933611
$ws.Range("O5").Value = 16.72783268763913
$ws.Range("B6").Value = 10.96279170436719
$ws.Range("C6").Value = 8.705192622442803
$ws.Range("E6").Value = 14.75385783998975
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.588637630859269
$ws.Range("I6").Value = 16.08385585296508
$ws.Range("N6").Value = 15.56118864256058
$ws.Range("O6").Value = 16.72958967999038
$ws.Range("B7").Value = 11.16818027006864
$ws.Range("C7").Value = 8.865251405484914
$ws.Range("E7").Value = 15.04715488214268
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.587959498010779
$ws.Range("I7").Value = 16.04929379728078
$ws.Range("N7").Value = 15.55017170709076
$ws.Range("O7").Value = 16.71794357880276
$ws.Range("B8").Value = 12.03296916709296
$ws.Range("C8").Value = 9.535266287276452
$ws.Range("E8").Value = 16.28683054451543
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.585118415246613
$ws.Range("I8").Value = 15.90725442750814
$ws.Range("N8").Value = 15.5055651738396
$ws.Range("O8").Value = 16.68021112280664
$ws.Range("B9").Value = 13.57012586331155
$ws.Range("C9").Value = 10.71463903908721
$ws.Range("E9").Value = 18.62900314657128
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.5800870359641
$ws.Range("I9").Value = 15.66668839580697
$ws.Range("N9").Value = 15.43245424763818
$ws.Range("O9").Value = 16.65585969934986
$ws.Range("B10").Value = 14.59668252658909
$ws.Range("C10").Value = 11.49657926549143
$ws.Range("E10").Value = 20.28671341568024
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.576715798791598
$ws.Range("I10").Value = 15.51339904627881
$ws.Range("N10").Value = 15.38745481184193
$ws.Range("O10").Value = 16.66863877579193
$ws.Range("B11").Value = 15.04048240294008
$ws.Range("C11").Value = 11.83359640740259
$ws.Range("E11").Value = 20.99834620173262
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.575251949781827
$ws.Range("I11").Value = 15.44882544197492
$ws.Range("N11").Value = 15.36886740162222
$ws.Range("O11").Value = 16.68117406443902
$ws.Range("B12").Value = 15.20966981783195
$ws.Range("C12").Value = 11.95850784116499
$ws.Range("E12").Value = 21.26175338260904
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.57470759437427
$ws.Range("I12").Value = 15.42512092434799
$ws.Range("N12").Value = 15.36209894979325
$ws.Range("O12").Value = 16.68689049899634
$ws.Range("B13").Value = 15.1711819926421
$ws.Range("C13").Value = 11.93172672858955
$ws.Range("E13").Value = 21.20529347765094
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.574824388422434
$ws.Range("I13").Value = 15.43019274171641
$ws.Range("N13").Value = 15.36354464916194
$ws.Range("O13").Value = 16.68561620631814
$ws.Range("B14").Value = 15.05409787559834
$ws.Range("C14").Value = 11.84392735637882
$ws.Range("E14").Value = 21.02013838051936
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.57520696579559
$ws.Range("I14").Value = 15.44686022710449
$ws.Range("N14").Value = 15.36830514519761
$ws.Range("O14").Value = 16.6816249145303
$ws.Range("B15").Value = 14.9827613866601
$ws.Range("C15").Value = 11.78979427985758
$ws.Range("E15").Value = 20.90593570009642
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.575442602514714
$ws.Range("I15").Value = 15.45716715897793
$ws.Range("N15").Value = 15.37125625683894
$ws.Range("O15").Value = 16.67930646839744
$ws.Range("B16").Value = 14.5672074939014
$ws.Range("C16").Value = 11.47417628281317
$ws.Range("E16").Value = 20.23935489847918
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.576812863342299
$ws.Range("I16").Value = 15.5177234215878
$ws.Range("N16").Value = 15.38870738468345
$ws.Range("O16").Value = 16.6679551137828
$ws.Range("B17").Value = 14.30629600336574
$ws.Range("C17").Value = 11.27575011117826
$ws.Range("E17").Value = 19.81957646389536
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.577671296634238
$ws.Range("I17").Value = 15.5561983835646
$ws.Range("N17").Value = 15.39989496538254
$ws.Range("O17").Value = 16.66271563012744
$ws.Range("B18").Value = 14.15404674658973
$ws.Range("C18").Value = 11.15986219726162
$ws.Range("E18").Value = 19.57413170298703
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.57817161272799
$ws.Range("I18").Value = 15.5788133167221
$ws.Range("N18").Value = 15.40650704283418
$ws.Range("O18").Value = 16.66033463103321
$ws.Range("B19").Value = 14.10212529480502
$ws.Range("C19").Value = 11.12032329165879
$ws.Range("E19").Value = 19.49034039559795
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.578342140946712
$ws.Range("I19").Value = 15.58655348021744
$ws.Range("N19").Value = 15.40877624405474
$ws.Range("O19").Value = 16.65963700551632
$ws.Range("B20").Value = 14.33429653285208
$ws.Range("C20").Value = 11.29705508324777
$ws.Range("E20").Value = 19.864676260453
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.577579235594593
$ws.Range("I20").Value = 15.55205240139576
$ws.Range("N20").Value = 15.39868568504965
$ws.Range("O20").Value = 16.66320788320108
$ws.Range("B21").Value = 15.08818563130435
$ws.Range("C21").Value = 11.86978984503988
$ws.Range("E21").Value = 21.07468746153428
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.575094323348051
$ws.Range("I21").Value = 15.44194422805212
$ws.Range("N21").Value = 15.36689954367863
$ws.Range("O21").Value = 16.68277092204989
$ws.Range("B22").Value = 15.59688974759206
$ws.Range("C22").Value = 12.22830494184159
$ws.Range("E22").Value = 21.83011900286271
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.573528390510349
$ws.Range("I22").Value = 15.37434575252875
$ws.Range("N22").Value = 15.34770011005638
$ws.Range("O22").Value = 16.70120847100127
$ws.Range("B23").Value = 15.32444405073597
$ws.Range("C23").Value = 12.03841018921522
$ws.Range("E23").Value = 21.43015709307209
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.574358860725897
$ws.Range("I23").Value = 15.41002293106569
$ws.Range("N23").Value = 15.3578033206665
$ws.Range("O23").Value = 16.69085020799722
$ws.Range("B24").Value = 14.32164449034618
$ws.Range("C24").Value = 11.2874287417499
$ws.Range("E24").Value = 19.84429944298684
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.577620835230173
$ws.Range("I24").Value = 15.55392525789461
$ws.Range("N24").Value = 15.39923183933031
$ws.Range("O24").Value = 16.66298336940789
$ws.Range("B25").Value = 13.17197918975378
$ws.Range("C25").Value = 10.41028946964173
$ws.Range("E25").Value = 17.98064983371976
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.581390750072044
$ws.Range("I25").Value = 15.72767171625578
$ws.Range("N25").Value = 15.45069931474101
$ws.Range("O25").Value = 16.65708621460956
